# "Generate Report for Handback" — mark the c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md
# row as handed back on the Overview sheet and on each per-language sheet, filling
# in the target file / handback file / handback datetime columns for that row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: row 7 is the c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md file.
# Columns E (zh-cn) and F (de-de) hold the per-language status text.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = $statusHandedBack
$overview.Range("F7").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: row 7 is the c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md file.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = $statusHandedBack
$zhcn.Range("I7").Value = "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md"
$zhcn.Range("J7").Value = "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.255b2d5e764e66770d8a2ff415421d3c53fa1eb1.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-11-09 00:46:03"
$zhcn.Hyperlinks.Add($zhcn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/364482c8c13f7817f35a89b7f0f23be6155faffc/e2e/c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md", "", "", "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md")

# ---------------------------------------------------------------------------
# de-de sheet: row 7 is the c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md file.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = $statusHandedBack
$dede.Range("I7").Value = "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md"
$dede.Range("J7").Value = "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.255b2d5e764e66770d8a2ff415421d3c53fa1eb1.de-de.xlf"
$dede.Range("K7").Value = "2016-11-09 00:46:21"
$dede.Hyperlinks.Add($dede.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/364482c8c13f7817f35a89b7f0f23be6155faffc/e2e/c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md", "", "", "c1a14f8c-ef6c-4fb3-90a7-fcb50d4918f7.md")
